# Update "想去人数" (want-to-go count) figures in F column across sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 7298
    5  = 28
    6  = 572
    7  = 202
    8  = 141
    11 = 64
    12 = 229
    13 = 19
    14 = 470
    18 = 50
    19 = 3833
    26 = 2516
    27 = 27
    28 = 330
    33 = 30
    35 = 2
    38 = 42
    39 = 1502
    40 = 176
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    2  = 7298
    5  = 28
    7  = 572
    8  = 202
    9  = 141
    12 = 64
    13 = 229
    14 = 19
    15 = 470
    19 = 50
    20 = 3833
    27 = 2516
    28 = 27
    29 = 330
    34 = 30
    36 = 2
    39 = 42
    40 = 1502
    41 = 176
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
